# Update Name of Algo
# Applies updated KNN imputation result values to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value = -12.721
$ws.Range("E3").Value = 16.736
$ws.Range("C21").Value = -12.499
$ws.Range("C23").Value = -12.594
$ws.Range("E24").Value = 16.717
$ws.Range("C25").Value = -12.37
$ws.Range("D27").Value = -8.484
$ws.Range("D31").Value = -8.327000000000002
$ws.Range("D39").Value = -7.747
$ws.Range("D48").Value = -7.475
$ws.Range("D51").Value = -8.401999999999999
$ws.Range("D52").Value = -7.522
$ws.Range("C53").Value = -11.523
$ws.Range("D55").Value = -8.065
$ws.Range("D56").Value = -8.228999999999999
$ws.Range("C57").Value = -13.401
$ws.Range("D57").Value = -8.537000000000001
$ws.Range("E57").Value = 16.679
$ws.Range("C59").Value = -13.063
$ws.Range("E61").Value = 16.628
$ws.Range("C69").Value = -11.118
$ws.Range("E70").Value = 17.568
$ws.Range("D73").Value = -8.004000000000001
$ws.Range("C79").Value = -12.013
$ws.Range("C83").Value = -13.169
$ws.Range("E86").Value = 16.597
$ws.Range("D89").Value = -6.702
$ws.Range("D90").Value = -7.441999999999998
$ws.Range("C93").Value = -11.511
$ws.Range("E98").Value = 16.421
$ws.Range("E100").Value = 16.725
$ws.Range("E102").Value = 16.49
